$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1, styled like the existing header cells (A1/B1)
$ws.Range("C1").Value = "calss"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data column values C2:C8
$ws.Range("C2").Value = "GDP กับ Set100 ✅"
$ws.Range("C3").Value = "GDP Set100 Export ✅"
$ws.Range("C4").Value = "GDP Export mech✅"
$ws.Range("C5").Value = "GDP Set100 Import ✅"
$ws.Range("C6").Value = "GDP กับ import ✅"
$ws.Range("C7").Value = "GDP predict GDP ✅"
$ws.Range("C8").Value = "GDP Set100 Import export✅"
